$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Sheets.Item("ALC")
$ws.Range("H9").Value = 135
$ws.Range("H17").Value = 840.44446
$ws.Range("J17").Value = 676.1539
$ws.Range("L17").Value = 2028.4617
$ws.Range("N17").Value = -2364.4617
$ws.Range("H70").Value = 3000
$ws.Range("I70").Value = 3000
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 9000
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -8730
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 3000
$ws.Range("I73").Value = 3000
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 9000
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -8064
$ws.Range("N73").ClearContents()
$ws.Range("H86").Value = 2000
$ws.Range("I86").Value = 2000
$ws.Range("K86").Value = 2000
$ws.Range("M86").Value = -877
$ws.Range("H88").Value = 5249.75
$ws.Range("J88").Value = 5333
$ws.Range("L88").Value = 5333
$ws.Range("N88").Value = -6145
$ws.Range("H89").Value = 2000
$ws.Range("I89").Value = 2000
$ws.Range("K89").Value = 10000
$ws.Range("M89").Value = -4384
$ws.Range("H91").Value = 5249.75
$ws.Range("J91").Value = 5333
$ws.Range("L91").Value = 5333
$ws.Range("N91").Value = -8141
$ws.Range("H98").Value = 2044.8334
$ws.Range("I98").Value = 2044.8334
$ws.Range("K98").Value = 2044.8334
$ws.Range("M98").Value = -546.8334
$ws.Range("H100").Value = 2333.3333
$ws.Range("I100").Value = 2000
$ws.Range("K100").Value = 2000
$ws.Range("M100").Value = -1459
$ws.Range("H106").Value = 20529.6
$ws.Range("I106").Value = 18792.783
$ws.Range("J106").Value = 40503
$ws.Range("K106").Value = 18792.783
$ws.Range("L106").Value = 40503
$ws.Range("M106").Value = -18161.783
$ws.Range("N106").Value = -41765
$ws.Range("H118").Value = 449
$ws.Range("I118").Value = 479.8
$ws.Range("J118").Value = 295
$ws.Range("K118").Value = 1439.4
$ws.Range("L118").Value = 885
$ws.Range("M118").Value = 217.5999999999999
$ws.Range("N118").Value = -4199
$ws.Range("H122").Value = 2044.8334
$ws.Range("I122").Value = 2044.8334
$ws.Range("K122").Value = 6134.5002
$ws.Range("M122").Value = -3684.5002
$ws.Range("H132").Value = 1893.1923
$ws.Range("I132").Value = 1402.8334
$ws.Range("J132").Value = 7777.5
$ws.Range("K132").Value = 4208.5002
$ws.Range("L132").Value = 23332.5
$ws.Range("M132").Value = -1678.5002
$ws.Range("N132").Value = -28392.5
$ws.Range("H137").Value = 2073.4
$ws.Range("I137").Value = 1973.5
$ws.Range("K137").Value = 5920.5
$ws.Range("M137").Value = -3370.5

# ---- Sheet: ARM ----
$ws = $wb.Sheets.Item("ARM")
$ws.Range("H61").Value = 4077.524
$ws.Range("I61").Value = 4121.4
$ws.Range("K61").Value = 4121.4
$ws.Range("M61").Value = -3909.4
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H102").Value = 2746.875
$ws.Range("I102").Value = 2193.75
$ws.Range("J102").Value = 3300
$ws.Range("K102").Value = 2193.75
$ws.Range("L102").Value = 3300
$ws.Range("M102").Value = -571.75
$ws.Range("N102").Value = -6544
$ws.Range("H136").Value = 4077.524
$ws.Range("I136").Value = 4121.4
$ws.Range("K136").Value = 12364.2
$ws.Range("M136").Value = -9814.199999999999

# ---- Sheet: BSM ----
$ws = $wb.Sheets.Item("BSM")
$ws.Range("H22").Value = 2720.75
$ws.Range("I22").Value = 294.33334
$ws.Range("J22").Value = 10000
$ws.Range("K22").Value = 294.33334
$ws.Range("L22").Value = 10000
$ws.Range("M22").Value = -121.33334
$ws.Range("N22").Value = -10346
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H86").Value = 2232.182
$ws.Range("I86").Value = 2232.182
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2232.182
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -1109.182
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 2232.182
$ws.Range("I89").Value = 2232.182
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 11160.91
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -5544.91
$ws.Range("N89").ClearContents()
$ws.Range("H99").Value = 3239.4
$ws.Range("I99").Value = 3239.4
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 3239.4
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1741.4
$ws.Range("N99").ClearContents()
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").ClearContents()
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Sheets.Item("CRP")
$ws.Range("H16").Value = 4266.4443
$ws.Range("J16").Value = 4050
$ws.Range("L16").Value = 4050
$ws.Range("N16").Value = -4624
$ws.Range("H22").Value = 9366.272000000001
$ws.Range("I22").Value = 280.8889
$ws.Range("J22").Value = 50250.5
$ws.Range("K22").Value = 280.8889
$ws.Range("L22").Value = 50250.5
$ws.Range("M22").Value = 69.11110000000002
$ws.Range("N22").Value = -50950.5
$ws.Range("H39").Value = 2000
$ws.Range("I39").Value = 2000
$ws.Range("K39").Value = 2000
$ws.Range("M39").Value = -1609
$ws.Range("H49").Value = 2000
$ws.Range("I49").Value = 2000
$ws.Range("K49").Value = 2000
$ws.Range("M49").Value = -1818
$ws.Range("H62").Value = 158603.6
$ws.Range("I62").Value = 98254.75
$ws.Range("K62").Value = 98254.75
$ws.Range("M62").Value = -97630.75
$ws.Range("H65").Value = 158603.6
$ws.Range("I65").Value = 98254.75
$ws.Range("K65").Value = 491273.75
$ws.Range("M65").Value = -488153.75
$ws.Range("H113").Value = 4266.4443
$ws.Range("J113").Value = 4050
$ws.Range("L113").Value = 4050
$ws.Range("N113").Value = -8390
$ws.Range("H134").Value = 2178.6365
$ws.Range("I134").Value = 1883.375
$ws.Range("J134").Value = 2966
$ws.Range("K134").Value = 5650.125
$ws.Range("L134").Value = 8898
$ws.Range("M134").Value = -3115.125
$ws.Range("N134").Value = -13968

# ---- Sheet: CUL ----
$ws = $wb.Sheets.Item("CUL")
$ws.Range("H12").Value = 147.77777
$ws.Range("J12").Value = 148.125
$ws.Range("L12").Value = 444.375
$ws.Range("N12").Value = -790.375
$ws.Range("H17").Value = 1866.6666
$ws.Range("J17").Value = 5000
$ws.Range("L17").Value = 15000
$ws.Range("N17").Value = -15338
$ws.Range("H23").Value = 250014
$ws.Range("I23").Value = 19
$ws.Range("J23").Value = 999999
$ws.Range("K23").Value = 57
$ws.Range("L23").Value = 2999997
$ws.Range("M23").Value = 178
$ws.Range("N23").Value = -3000467
$ws.Range("H24").Value = 1883.2
$ws.Range("I24").Value = 237.5
$ws.Range("J24").Value = 2294.625
$ws.Range("K24").Value = 712.5
$ws.Range("L24").Value = 6883.875
$ws.Range("M24").Value = -482.5
$ws.Range("N24").Value = -7343.875
$ws.Range("H34").Value = 21116.8
$ws.Range("J34").Value = 41660.6
$ws.Range("L34").Value = 124981.8
$ws.Range("N34").Value = -125149.8
$ws.Range("H55").Value = 53525
$ws.Range("J55").Value = 71333.336
$ws.Range("L55").Value = 214000.008
$ws.Range("N55").Value = -214354.008
$ws.Range("H68").Value = 1001.4
$ws.Range("I68").Value = 1001.3333
$ws.Range("K68").Value = 3003.9999
$ws.Range("M68").Value = -2192.9999
$ws.Range("H71").Value = 1001.4
$ws.Range("I71").Value = 1001.3333
$ws.Range("K71").Value = 9011.9997
$ws.Range("M71").Value = -4955.9997
$ws.Range("H122").Value = 338.4
$ws.Range("J122").Value = 364.66666
$ws.Range("L122").Value = 3281.99994
$ws.Range("N122").Value = -8181.99994
$ws.Range("H131").Value = 1474.8
$ws.Range("I131").Value = 907.5
$ws.Range("J131").Value = 1519.2941
$ws.Range("K131").Value = 2722.5
$ws.Range("L131").Value = 4557.8823
$ws.Range("M131").Value = 2317.5
$ws.Range("N131").Value = -14637.8823

# ---- Sheet: GSM ----
$ws = $wb.Sheets.Item("GSM")
$ws.Range("H122").Value = 57777.055
$ws.Range("I122").Value = 1853.3846
$ws.Range("J122").Value = 203178.6
$ws.Range("K122").Value = 5560.1538
$ws.Range("L122").Value = 609535.8
$ws.Range("M122").Value = -3110.1538
$ws.Range("N122").Value = -614435.8

# ---- Sheet: LTW ----
$ws = $wb.Sheets.Item("LTW")
$ws.Range("H22").Value = 1166.5714
$ws.Range("I22").Value = 1211
$ws.Range("K22").Value = 1211
$ws.Range("M22").Value = -916
$ws.Range("H27").Value = 1166.5714
$ws.Range("I27").Value = 1211
$ws.Range("K27").Value = 1211
$ws.Range("M27").Value = -1104
$ws.Range("H46").Value = 2116.1667
$ws.Range("I46").Value = 1633.3334
$ws.Range("J46").Value = 2277.111
$ws.Range("K46").Value = 1633.3334
$ws.Range("L46").Value = 2277.111
$ws.Range("M46").Value = -1445.3334
$ws.Range("N46").Value = -2653.111
$ws.Range("H55").Value = 814.5909
$ws.Range("I55").Value = 675.8461
$ws.Range("K55").Value = 675.8461
$ws.Range("M55").Value = -502.8461
$ws.Range("H136").Value = 6304.625
$ws.Range("I136").Value = 6283.857
$ws.Range("J136").Value = 6450
$ws.Range("K136").Value = 18851.571
$ws.Range("L136").Value = 19350
$ws.Range("M136").Value = -16301.571
$ws.Range("N136").Value = -24450

# ---- Sheet: WVR ----
$ws = $wb.Sheets.Item("WVR")
$ws.Range("H5").Value = 3000000
$ws.Range("I5").Value = 3000000
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 3000000
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -2999888
$ws.Range("N5").ClearContents()
